$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 118, shifting rows 118:224 down to 119:225.
$ws.Rows.Item(118).Insert()

# Populate the newly inserted row 118 with the new record's data.
$ws.Cells.Item(118, 1).Value = 10
$ws.Cells.Item(118, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(118, 3).Value = "La Araucanía"
$ws.Cells.Item(118, 4).Value = 44586
$ws.Cells.Item(118, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(118, 5).Value = 9
$ws.Cells.Item(118, 6).Value = 100112001
$ws.Cells.Item(118, 7).Value = "Berenjena"
$ws.Cells.Item(118, 8).Value = "Sin especificar"
$ws.Cells.Item(118, 9).Value = "Primera"
$ws.Cells.Item(118, 10).Value = 35
$ws.Cells.Item(118, 11).Value = 12000
$ws.Cells.Item(118, 12).Value = 12000
$ws.Cells.Item(118, 13).Value = 12000
$ws.Cells.Item(118, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(118, 15).Value = "Región del Maule"
$ws.Cells.Item(118, 16).Value = 200
$ws.Cells.Item(118, 17).Value = 60
$ws.Cells.Item(118, 18).Value = "Hortaliza"
